$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remove the old "average" / "total" summary rows (117 and 119) from Sheet1.
$ws1.Range("A117:J119").Delete()

# Update sheet1 view: deselect as active tab, move to bottom, select header row E1:J1
$ws1.Range("E1:J1").Select()

# Add a new worksheet after Sheet1 for the per-column averages.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "averages"

# Header row, reusing the same shared strings as Sheet1's column headers.
$ws2.Range("A1").Value = "zero_shot_time"
$ws2.Range("B1").Value = "tf_idf_time"
$ws2.Range("C1").Value = "confidence_score_time"
$ws2.Range("D1").Value = "text_generation_time"
$ws2.Range("E1").Value = "total_time"
$ws2.Range("F1").Value = "Bleu"
$ws2.Range("G1").Value = "Perplex"

# Data row with averages pulled from Sheet1, plus a total-time sum.
$ws2.Range("A2").Formula = "=AVERAGE(Sheet1!E2:E116)"
$ws2.Range("B2").Formula = "=AVERAGE(Sheet1!F2:F116)"
$ws2.Range("C2").Formula = "=AVERAGE(Sheet1!G2:G116)"
$ws2.Range("D2").Formula = "=AVERAGE(Sheet1!H2:H116)"
$ws2.Range("E2").Formula = "=SUM(A2:D2)"
$ws2.Range("F2").Formula = "=AVERAGE(Sheet1!I2:I116)"
$ws2.Range("G2").Formula = "=AVERAGE(Sheet1!J2:J116)"

$ws2.Range("E2").Select()
